$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last-updated timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 19:09"

# Row 4
$ws.Range("B4").Value = 8353795
$ws.Range("C4").Value = 11130
$ws.Range("D4").Value = 5439043
$ws.Range("E4").Value = 2690336
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 224416

# Row 5
$ws.Range("B5").Value = 7544353
$ws.Range("C5").Value = 51626
$ws.Range("D5").Value = 6654464
$ws.Range("E5").Value = 775314
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 511
$ws.Range("H5").Value = 114575

# Row 6
$ws.Range("B6").Value = 5224821
$ws.Range("C6").Value = 459
$ws.Range("D6").Value = 4635315
$ws.Range("E6").Value = 435776
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 153730

# Row 21
$ws.Range("B21").Value = 364303
$ws.Range("C21").Value = 2570
$ws.Range("D21").Value = 290000
$ws.Range("E21").Value = 64443
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 9860

# Row 24
$ws.Range("B24").Value = 347493
$ws.Range("C24").Value = 1815
$ws.Range("D24").Value = 304003
$ws.Range("E24").Value = 34194
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 9296

# Row 27
$ws.Range("B27").Value = 303109
$ws.Range("C27").Value = 339
$ws.Range("D27").Value = 268093
$ws.Range("E27").Value = 32807
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 19
$ws.Range("H27").Value = 2209

# Row 34
$ws.Range("A34").Value = "Chequia"
$ws.Range("B34").Value = 171487
$ws.Range("C34").Value = 2660
$ws.Range("D34").Value = 69090
$ws.Range("E34").Value = 100995
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 50
$ws.Range("H34").Value = 1402

# Row 35
$ws.Range("A35").Value = "Marruecos"
$ws.Range("B35").Value = 170911
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 141381
$ws.Range("E35").Value = 26652
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 2878

# Row 66
$ws.Range("B66").Value = 54402
$ws.Range("C66").Value = 199
$ws.Range("D66").Value = 38088
$ws.Range("E66").Value = 14458
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 10
$ws.Range("H66").Value = 1856

# Row 69
$ws.Range("A69").Value = "Irlanda"
$ws.Range("B69").Value = 49962
$ws.Range("C69").Value = 1283
$ws.Range("D69").Value = 23364
$ws.Range("E69").Value = 24746
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 1852

# Row 70
$ws.Range("A70").Value = "Libia"
$ws.Range("B70").Value = 48790
$ws.Range("C70").Value = 945
$ws.Range("D70").Value = 26889
$ws.Range("E70").Value = 21176
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 26
$ws.Range("H70").Value = 725

# Row 99
$ws.Range("B99").Value = 15615
$ws.Range("C99").Value = 86
$ws.Range("D99").Value = 10994
$ws.Range("E99").Value = 4385
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 5
$ws.Range("H99").Value = 236

# Row 108
$ws.Range("A108").Value = "Luxemburgo"
$ws.Range("B108").Value = 10888
$ws.Range("C108").Value = 242
$ws.Range("D108").Value = 8468
$ws.Range("E108").Value = 2287
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 133

# Row 109
$ws.Range("A109").Value = "Mozambique"
$ws.Range("B109").Value = 10866
$ws.Range("C109").Value = 159
$ws.Range("D109").Value = 8513
$ws.Range("E109").Value = 2278
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 75

# Row 124
$ws.Range("B124").Value = 5780
$ws.Range("C124").Value = 15
$ws.Range("D124").Value = 5415
$ws.Range("E124").Value = 249
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 116

# Row 126
$ws.Range("B126").Value = 5538
$ws.Range("C126").Value = 63
$ws.Range("D126").Value = 3403
$ws.Range("E126").Value = 2122
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 13

# Row 154
$ws.Range("B154").Value = 2644
$ws.Range("C154").Value = 63
$ws.Range("D154").Value = 1444
$ws.Range("E154").Value = 1175
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 25
